# Generate Report for Handback
# Update timestamps for the 553ec807 (handback) file rows across sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the 553ec807 file row (row 2)
$wsOverview.Range("G2").Value = "2016-08-21 20:59:16"

# zh-cn sheet: row 2 is the 553ec807 file
$wsZhCn.Range("H2").Value = "2016-08-21 20:59:12"
$wsZhCn.Range("K2").Value = "2016-08-21 20:59:28"

# de-de sheet: row 2 is the 553ec807 file
$wsDeDe.Range("H2").Value = "2016-08-21 20:59:16"
$wsDeDe.Range("K2").Value = "2016-08-21 20:59:34"
